# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the data refresh captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Cell -> new value map (same updates applied identically to both sheets)
$updates = @{
    "F2"  = 4918
    "F3"  = 136
    "F4"  = 110
    "F5"  = 810
    "F6"  = 249
    "F7"  = 1289
    "F8"  = 139
    "F15" = 4364
    "F16" = 6641
    "F18" = 50
    "F21" = 52
    "F22" = 4077
    "F23" = 438
    "F24" = 63
    "F25" = 42
    "F26" = 2661
    "F27" = 573
    "F29" = 158
    "F30" = 331
    "F31" = 341
    "F33" = 208
    "F34" = 29
    "F35" = 1603
    "F36" = 1004
    "F38" = 115
    "F42" = 9
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
